# Generate Report for Handoff
# Rotates the handoff artifact GUID from 9009084b-...-432457 to
# d0ca1dd6-...-cd4d5, refreshes the handoff hash + timestamps.

$wb = $excel.ActiveWorkbook

$oldGuid = "9009084b-b938-466d-8fe2-db8bcc432457"
$newGuid = "d0ca1dd6-2ade-4e81-a460-1e79887cd4d5"
$oldHash = "905381bf0ae7f3a162d3e46b42f7788d05f0310e"
$newHash = "4f97ae72bdc331ac2cf7ad237d892bf152f05406"

# NOTE: the original commit left the hyperlink *target* (relationship URL,
# still pinned to the old commit SHA / old GUID) untouched and only
# refreshed the cell text + the hyperlink's visible "display" text - so we
# reuse the existing (old-GUID) target URL rather than re-pointing it at
# the new GUID.
$repoBlobBase = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/3c358fbba47d84fc14192691056b4d4525a67425/e2e/"
$origTarget = "$repoBlobBase$oldGuid.md"

# --- Overview sheet ---
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("A2").Value = "$newGuid.md"
$wsOverview.Range("B2").Value = "e2e\$newGuid.md"
$wsOverview.Range("G2").Value = "2016-08-22 16:59:11"

$wsOverview.Hyperlinks.Delete()
$wsOverview.Hyperlinks.Add($wsOverview.Range("B2"), $origTarget, "", "", "e2e\$newGuid.md")

# --- zh-cn sheet ---
$wsZh = $wb.Worksheets.Item("zh-cn")
$wsZh.Range("A2").Value = "$newGuid.md"
$wsZh.Range("G2").Value = "$newGuid.$newHash.zh-cn.xlf"
$wsZh.Range("H2").Value = "2016-08-22 16:58:59"

$wsZh.Hyperlinks.Delete()
$wsZh.Hyperlinks.Add($wsZh.Range("A2"), $origTarget, "", "", "$newGuid.md")

# --- de-de sheet ---
$wsDe = $wb.Worksheets.Item("de-de")
$wsDe.Range("A2").Value = "$newGuid.md"
$wsDe.Range("G2").Value = "$newGuid.$newHash.de-de.xlf"
$wsDe.Range("H2").Value = "2016-08-22 16:59:11"

$wsDe.Hyperlinks.Delete()
$wsDe.Hyperlinks.Add($wsDe.Range("A2"), $origTarget, "", "", "$newGuid.md")
